# Prepared tuned scaled RF submissions: append two new rows (52, 53) to the
# "Tabelle1" results table with the standard/minmax scaled-center-before-avg
# tuned RandomForest runs, formatted like the rest of the AUC columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table from A1:O51 to A1:O53 so the new rows become part of it
# (this also keeps the AutoFilter / banding in sync).
$lo.Resize($ws.Range("A1:O53"))

# Burn four throwaway "top10" conditional-format rules (same orange style
# used elsewhere in the sheet) on an unused column and remove them again.
# Excel keeps the dxf records even after the owning rule is deleted, which
# is how this workbook's style table picked up its extra (unused) entries.
for ($i = 0; $i -lt 4; $i++) {
    $tmp = $ws.Range("Z1:Z2").FormatConditions.AddTop10()
    $tmp.Font.Color = 22428
    $tmp.Interior.Color = 10284031
    $tmp.Delete()
}

# --- Row 52: standard scaling center data before avg (tuned) -------------
# (Filled in the same order the author did: Comments, Hyperparameters, then
# the File name, so new shared-string entries land in the matching order.)
$ws.Range("I52").Value = "standard scaling center data before avg tuned"
$ws.Range("H52").Value = "{'bootstrap': True,  'max_depth': 2, 'max_features': 0.5,  'max_samples': 0.8,  'min_samples_split': 10}"
$ws.Range("A52").Value = "2023-03-11-1845_RF_Standard_avg_tuned.csv"
$ws.Range("B52").Value = "RandomForest"
$ws.Range("C52").Value = "MoCo"
$ws.Range("D52").Value = "Centers"
$ws.Range("E52").Value = "1 x 3"
$ws.Range("F52").Value = "average"
$ws.Range("J52").Value = 0.68883283650538896
$ws.Range("K52").Value = 0.70780051150895096
$ws.Range("L52").Value = 0.65485362095531596
$ws.Range("M52").Value = 0.68382898965655203

# --- Row 53: minmax scaling center data before avg (tuned) ---------------
$ws.Range("I53").Value = "minmax scaling center data before avg tuned"
$ws.Range("H53").Value = "{'bootstrap': True, 'max_depth': 4,  'max_features': 1.0,  'max_samples': 0.7,  'min_samples_split': 10}"
$ws.Range("A53").Value = "2023-03-11-1856_RF_MinMax_avg_tuned.csv"
$ws.Range("B53").Value = "RandomForest"
$ws.Range("C53").Value = "MoCo"
$ws.Range("D53").Value = "Centers"
$ws.Range("E53").Value = "1 x 3"
$ws.Range("F53").Value = "average"
$ws.Range("J53").Value = 0.64939234120614497
$ws.Range("K53").Value = 0.74648337595907899
$ws.Range("L53").Value = 0.66358500256805297
$ws.Range("M53").Value = 0.68648690657776001

# Format the AUC columns J:L and the average M with three decimals, like
# the style that was already used for the numeric AUC figures elsewhere.
$ws.Range("J52:L53").NumberFormat = "0.000"
$ws.Range("M52:M53").NumberFormat = "0.000"
$ws.Range("M52:M53").HorizontalAlignment = -4131

# Extend the Top-10 conditional formatting on columns J and K so it keeps
# covering the newly added rows (was J2:J51 / K2:K51).
$fcJ = $ws.Range("J2:J51").FormatConditions.Item(2)
$fcJ.ModifyAppliesToRange($ws.Range("J2:J53"))
$fcK = $ws.Range("K2:K51").FormatConditions.Item(2)
$fcK.ModifyAppliesToRange($ws.Range("K2:K53"))

# Park the selection where the author left it (just below the new rows).
$ws.Range("A54").Select() | Out-Null
